$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 595
$ws.Range("A595").Value = 'outputs/2024-04-30/06-08-30'
$ws.Range("B595").Value = $false
$ws.Range("C595").Value = 'eicu'
$ws.Range("D595").Value = 'diagnosis'
$ws.Range("E595").Value = 'VA'
$ws.Range("F595").Value = 'descemb_bert'
$ws.Range("G595").Value = 'ehr_model'
$ws.Range("H595").Value = "'False"
$ws.Range("I595").Value = "'False"
$ws.Range("J595").Value = "'False"

# Row 596
$ws.Range("A596").Value = 'outputs/2024-04-30/06-08-33'
$ws.Range("B596").Value = $false
$ws.Range("C596").Value = 'eicu'
$ws.Range("D596").Value = 'diagnosis'
$ws.Range("E596").Value = 'DSVA'
$ws.Range("F596").Value = 'descemb_bert'
$ws.Range("G596").Value = 'ehr_model'
$ws.Range("H596").Value = "'False"
$ws.Range("I596").Value = "'False"
$ws.Range("J596").Value = "'False"

# Row 597
$ws.Range("A597").Value = 'outputs/2024-04-30/06-12-43'
$ws.Range("B597").Value = $false
$ws.Range("C597").Value = 'eicu'
$ws.Range("D597").Value = 'diagnosis'
$ws.Range("E597").Value = 'VA'
$ws.Range("F597").Value = 'descemb_bert'
$ws.Range("G597").Value = 'ehr_model'
$ws.Range("H597").Value = "'False"
$ws.Range("I597").Value = "'False"
$ws.Range("J597").Value = "'False"

# Row 598
$ws.Range("A598").Value = 'outputs/2024-04-30/06-14-20'
$ws.Range("B598").Value = $true
$ws.Range("C598").Value = 'eicu'
$ws.Range("D598").Value = 'diagnosis'
$ws.Range("E598").Value = 'VA'
$ws.Range("F598").Value = 'descemb_bert'
$ws.Range("G598").Value = 'ehr_model'
$ws.Range("H598").Value = "'False"
$ws.Range("I598").Value = "'False"
$ws.Range("J598").Value = "'False"
$ws.Range("K598").Value = 241
$ws.Range("L598").Value = 0.823
$ws.Range("M598").Value = 0.772
$ws.Range("N598").Value = 0.443

# Row 599
$ws.Range("A599").Value = 'outputs/2024-04-30/07-42-04'
$ws.Range("B599").Value = $true
$ws.Range("C599").Value = 'eicu'
$ws.Range("D599").Value = 'diagnosis'
$ws.Range("E599").Value = 'DSVA'
$ws.Range("F599").Value = 'descemb_bert'
$ws.Range("G599").Value = 'ehr_model'
$ws.Range("H599").Value = "'False"
$ws.Range("I599").Value = "'False"
$ws.Range("J599").Value = "'False"
$ws.Range("K599").Value = 224
$ws.Range("L599").Value = 0.725
$ws.Range("M599").Value = 0.802
$ws.Range("N599").Value = 0.484

# Row 600
$ws.Range("A600").Value = 'outputs/2024-04-30/09-06-18'
$ws.Range("B600").Value = $false
$ws.Range("C600").Value = 'eicu'
$ws.Range("D600").Value = 'diagnosis'
$ws.Range("E600").Value = 'DSVA_DPE'
$ws.Range("F600").Value = 'descemb_bert'
$ws.Range("G600").Value = 'ehr_model'
$ws.Range("H600").Value = "'False"
$ws.Range("I600").Value = "'False"
$ws.Range("J600").Value = "'False"

# Row 601
$ws.Range("A601").Value = 'outputs/2024-04-30/09-06-21'
$ws.Range("B601").Value = $true
$ws.Range("C601").Value = 'eicu'
$ws.Range("D601").Value = 'diagnosis'
$ws.Range("E601").Value = 'VC'
$ws.Range("F601").Value = 'descemb_bert'
$ws.Range("G601").Value = 'ehr_model'
$ws.Range("H601").Value = "'False"
$ws.Range("I601").Value = "'False"
$ws.Range("J601").Value = "'False"
$ws.Range("K601").Value = 270
$ws.Range("L601").Value = 0.6830000000000001
$ws.Range("M601").Value = 0.785
$ws.Range("N601").Value = 0.447

# Row 602
$ws.Range("A602").Value = 'outputs/2024-04-30/10-37-13'
$ws.Range("B602").Value = $true
$ws.Range("C602").Value = 'eicu'
$ws.Range("D602").Value = 'los_3day'
$ws.Range("E602").Value = 'VA'
$ws.Range("F602").Value = 'descemb_bert'
$ws.Range("G602").Value = 'ehr_model'
$ws.Range("H602").Value = "'False"
$ws.Range("I602").Value = "'False"
$ws.Range("J602").Value = "'False"
$ws.Range("K602").Value = 244
$ws.Range("L602").Value = 3.399
$ws.Range("M602").Value = 0.698
$ws.Range("N602").Value = 0.496

# Row 603
$ws.Range("A603").Value = 'outputs/2024-04-30/12-03-53'
$ws.Range("B603").Value = $true
$ws.Range("C603").Value = 'eicu'
$ws.Range("D603").Value = 'los_3day'
$ws.Range("E603").Value = 'DSVA'
$ws.Range("F603").Value = 'descemb_bert'
$ws.Range("G603").Value = 'ehr_model'
$ws.Range("H603").Value = "'False"
$ws.Range("I603").Value = "'False"
$ws.Range("J603").Value = "'False"
$ws.Range("K603").Value = 103
$ws.Range("L603").Value = 2.338
$ws.Range("M603").Value = 0.695
$ws.Range("N603").Value = 0.481

# Row 604
$ws.Range("A604").Value = 'outputs/2024-04-30/12-41-35'
$ws.Range("B604").Value = $false
$ws.Range("C604").Value = 'eicu'
$ws.Range("D604").Value = 'los_3day'
$ws.Range("E604").Value = 'DSVA_DPE'
$ws.Range("F604").Value = 'descemb_bert'
$ws.Range("G604").Value = 'ehr_model'
$ws.Range("H604").Value = "'False"
$ws.Range("I604").Value = "'False"
$ws.Range("J604").Value = "'False"

# Row 605
$ws.Range("A605").Value = 'outputs/2024-04-30/12-41-38'
$ws.Range("B605").Value = $true
$ws.Range("C605").Value = 'eicu'
$ws.Range("D605").Value = 'los_3day'
$ws.Range("E605").Value = 'VC'
$ws.Range("F605").Value = 'descemb_bert'
$ws.Range("G605").Value = 'ehr_model'
$ws.Range("H605").Value = "'False"
$ws.Range("I605").Value = "'False"
$ws.Range("J605").Value = "'False"
$ws.Range("K605").Value = 128
$ws.Range("L605").Value = 2.395
$ws.Range("M605").Value = 0.6919999999999999
$ws.Range("N605").Value = 0.465

# Row 606
$ws.Range("A606").Value = 'outputs/2024-04-30/13-23-34'
$ws.Range("B606").Value = $true
$ws.Range("C606").Value = 'eicu'
$ws.Range("D606").Value = 'los_7day'
$ws.Range("E606").Value = 'VA'
$ws.Range("F606").Value = 'descemb_bert'
$ws.Range("G606").Value = 'ehr_model'
$ws.Range("H606").Value = "'False"
$ws.Range("I606").Value = "'False"
$ws.Range("J606").Value = "'False"
$ws.Range("K606").Value = 244
$ws.Range("L606").Value = 2.006
$ws.Range("M606").Value = 0.641
$ws.Range("N606").Value = 0.165

# Row 607
$ws.Range("A607").Value = 'outputs/2024-04-30/14-50-12'
$ws.Range("B607").Value = $true
$ws.Range("C607").Value = 'eicu'
$ws.Range("D607").Value = 'los_7day'
$ws.Range("E607").Value = 'DSVA'
$ws.Range("F607").Value = 'descemb_bert'
$ws.Range("G607").Value = 'ehr_model'
$ws.Range("H607").Value = "'False"
$ws.Range("I607").Value = "'False"
$ws.Range("J607").Value = "'False"
$ws.Range("K607").Value = 151
$ws.Range("L607").Value = 1.758
$ws.Range("M607").Value = 0.649
$ws.Range("N607").Value = 0.168

# Row 608
$ws.Range("A608").Value = 'outputs/2024-04-30/15-45-20'
$ws.Range("B608").Value = $false
$ws.Range("C608").Value = 'eicu'
$ws.Range("D608").Value = 'los_7day'
$ws.Range("E608").Value = 'DSVA_DPE'
$ws.Range("F608").Value = 'descemb_bert'
$ws.Range("G608").Value = 'ehr_model'
$ws.Range("H608").Value = "'False"
$ws.Range("I608").Value = "'False"
$ws.Range("J608").Value = "'False"

# Row 609
$ws.Range("A609").Value = 'outputs/2024-04-30/15-45-22'
$ws.Range("B609").Value = $false
$ws.Range("C609").Value = 'eicu'
$ws.Range("D609").Value = 'los_7day'
$ws.Range("E609").Value = 'VC'
$ws.Range("F609").Value = 'descemb_bert'
$ws.Range("G609").Value = 'ehr_model'
$ws.Range("H609").Value = "'False"
$ws.Range("I609").Value = "'False"
$ws.Range("J609").Value = "'False"
$ws.Range("K609").Value = 60
$ws.Range("L609").Value = 0.278
$ws.Range("M609").Value = 0.87
$ws.Range("N609").Value = 0.62
